$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Locate the two target paragraphs by scanning the Paragraphs collection
# (more reliable in this host than Find + Paragraphs.Item(1).Next()).
# ---------------------------------------------------------------------------
$total = $d.Paragraphs.Count
$idxLastYG = -1        # paragraph ending in "...放在最后一个。"
$idxJuanZhanLan = -1   # paragraph ending in "...卷展栏"
for ($i = 1; $i -le $total; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*最后一个。*") {
        $idxLastYG = $i
    }
    if ($t -like "*卷展栏*") {
        $idxJuanZhanLan = $i
    }
}

# ---------------------------------------------------------------------------
# Change 1: simplify the empty paragraph right after "...最后一个。" from a
# paragraph carrying an empty <w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr>
# down to a bare empty paragraph.
# ---------------------------------------------------------------------------
if ($idxLastYG -gt 0) {
    $emptyIdx = $idxLastYG + 1
    $emptyRng = $d.Paragraphs.Item($emptyIdx).Range

    $blankXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    $emptyRng.InsertXML($blankXml) | Out-Null
}

# ---------------------------------------------------------------------------
# Change 2: after the paragraph that ends with "卷展栏" (and before the final
# bookmark paragraph), add a new "153 EditAnywhere vs EditDefaultsOnly"
# heading plus two body paragraphs describing the two specifiers.
# ---------------------------------------------------------------------------
if ($idxJuanZhanLan -gt 0) {
    $anchor = $d.Paragraphs.Item($idxJuanZhanLan)

    # Create three brand-new empty paragraphs right after the anchor
    # paragraph (this correctly splits a new paragraph mark out of the
    # anchor instead of merging into the following paragraph).
    $anchor.Range.InsertParagraphAfter() | Out-Null
    $anchor.Range.InsertParagraphAfter() | Out-Null
    $anchor.Range.InsertParagraphAfter() | Out-Null

    $headingIdx = $idxJuanZhanLan + 1
    $para2Idx = $idxJuanZhanLan + 2
    $para3Idx = $idxJuanZhanLan + 3

    # --- heading paragraph: "153 EditAnywhere vs EditDefaultsOnly" ---
    $headingXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>153 EditAnywhere vs EditDefaultsOnly</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $d.Paragraphs.Item($headingIdx).Range.InsertXML($headingXml) | Out-Null
    $d.Paragraphs.Item($headingIdx).Range.Style = "1"

    # --- body paragraph: EditAnywhere description ---
    $para2Xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>Edit</w:t></w:r><w:r><w:t xml:space="preserve">Anywhere </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>是每个实例都可以修改单独的数据</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $d.Paragraphs.Item($para2Idx).Range.InsertXML($para2Xml) | Out-Null

    # --- body paragraph: EditDefaultsOnly description ---
    $para3Xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>EditDefaultsOnly</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>只能在蓝图中修改默认值，也就是所有实例必须保持一致</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $d.Paragraphs.Item($para3Idx).Range.InsertXML($para3Xml) | Out-Null
}
